$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output3Salted")

$values = @(1473,4124,2604,-1555,4909,-1476,4505,5256,-2839,-764,-2136,-2245,-2660,4034,-154,-2506,-840,-947,5370,28,-393,1117,-1244,7352,1489,8707,7433,871,2522,3367,5348,10364,12427,3973,11106,6837,6632,7341,13990,14910,14916,9940,16891,18179,16685,18818,17420,19446,18201,12922,12931,19945,23006,16328,22620,22685,26678,21718,22350,28775,23064,21565,23344,32651,25739,26439,27741,36676,30966,32875,37387,38522,35864,39299,34075,42025,38064,38146,38294,39664,44088,43976,49176,46689,45226,48625,49137,57303,51966,52887,55568,62349,64010,59495,59115,59384,62365,64474,64510,72414)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}
